$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 3, pushing the existing data (rows 3-6) down to rows 4-7
$ws.Rows.Item(3).Insert()

# Select the newly inserted (now empty) row 3, matching the post-edit selection state
$ws.Rows.Item(3).Select()
